$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value = "ano"
$ws.Range("L7").Value = "GV pro Stunde"

$ws.Range("Q6").Select()
